$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new row above row 31 (pushes existing row 31 "danholland" and
# everything below it down by one row, to row 32 and beyond).
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new default production user.
$ws.Cells.Item(31, 1).Value = "Oleg_Babak"
$ws.Cells.Item(31, 2).Value = "Password1!"

# Reflect the cell that the user ended up selecting after the edit.
$ws.Range("F31").Select()
